# Update iqma (congregation) prayer times: new isha and fajr times.
# Source data: diff of xl/worksheets/sheet1.xml — column L (ishaIqma) moves
# from ~19:00/19:15 to 18:30, and column C (fajrIqma) moves from 6:45 to
# 6:30 for the affected rows (rows 2-32 and 325-367).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newTimes = @{
    "L2" = 0.77083333333333304
    "L3" = 0.77083333333333337
    "L4" = 0.77083333333333404
    "L5" = 0.77083333333333404
    "L6" = 0.77083333333333404
    "L7" = 0.77083333333333504
    "L8" = 0.77083333333333504
    "L9" = 0.77083333333333504
    "L10" = 0.77083333333333603
    "L11" = 0.77083333333333603
    "L12" = 0.77083333333333603
    "L13" = 0.77083333333333703
    "L14" = 0.77083333333333703
    "L15" = 0.77083333333333703
    "L16" = 0.77083333333333803
    "L17" = 0.77083333333333803
    "L18" = 0.77083333333333803
    "L19" = 0.77083333333333903
    "L20" = 0.77083333333333903
    "L21" = 0.77083333333333903
    "L22" = 0.77083333333334003
    "L23" = 0.77083333333334003
    "L24" = 0.77083333333334003
    "L25" = 0.77083333333334103
    "L26" = 0.77083333333334103
    "L27" = 0.77083333333334103
    "L28" = 0.77083333333334203
    "L29" = 0.77083333333334203
    "L30" = 0.77083333333334203
    "L31" = 0.77083333333334303
    "L32" = 0.77083333333334303
    "L325" = 0.77083333333333337
    "C326" = 0.27083333333333298
    "L326" = 0.77083333333333337
    "C327" = 0.27083333333333298
    "L327" = 0.77083333333333304
    "C328" = 0.27083333333333298
    "L328" = 0.77083333333333304
    "C329" = 0.27083333333333298
    "L329" = 0.77083333333333304
    "C330" = 0.27083333333333298
    "L330" = 0.77083333333333304
    "C331" = 0.27083333333333298
    "L331" = 0.77083333333333304
    "L332" = 0.77083333333333304
    "L333" = 0.77083333333333304
    "L334" = 0.77083333333333304
    "L335" = 0.77083333333333304
    "L336" = 0.77083333333333304
    "L337" = 0.77083333333333304
    "L338" = 0.77083333333333304
    "L339" = 0.77083333333333304
    "L340" = 0.77083333333333304
    "L341" = 0.77083333333333304
    "L342" = 0.77083333333333304
    "L343" = 0.77083333333333304
    "L344" = 0.77083333333333304
    "L345" = 0.77083333333333304
    "L346" = 0.77083333333333304
    "L347" = 0.77083333333333304
    "L348" = 0.77083333333333304
    "L349" = 0.77083333333333304
    "L350" = 0.77083333333333304
    "L351" = 0.77083333333333304
    "L352" = 0.77083333333333304
    "L353" = 0.77083333333333304
    "L354" = 0.77083333333333304
    "L355" = 0.77083333333333304
    "L356" = 0.77083333333333304
    "L357" = 0.77083333333333304
    "L358" = 0.77083333333333304
    "L359" = 0.77083333333333304
    "L360" = 0.77083333333333304
    "L361" = 0.77083333333333304
    "L362" = 0.77083333333333304
    "L363" = 0.77083333333333304
    "L364" = 0.77083333333333304
    "L365" = 0.77083333333333304
    "L366" = 0.77083333333333304
    "L367" = 0.77083333333333337
}

foreach ($cellRef in $newTimes.Keys) {
    $ws.Range($cellRef).Value = $newTimes[$cellRef]
}

# Restore the author's last on-screen selection/scroll position for Sheet1.
$ws.Activate()
$ws.Range("M375").Select()
$excel.ActiveWindow.ScrollRow = 319
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Updated $($newTimes.Count) iqma time cells (isha/fajr)."
